# Refresh the cryptos price/volume table (and fix the Aptos/EnergySwap row
# order) per the GitHub Actions data-refresh commit.
#
# All data cells on this sheet are plain text (coinranking.com renders
# prices like "63.292.13" and volumes like "  +2.70%  " as strings, not
# numbers). Excel's smart cell-entry will happily reinterpret a
# numeric-looking string (e.g. "584.66", "7.70", "0.0770") as a real
# number and silently drop formatting such as trailing zeros, so any
# value that parses as a number is written with a leading apostrophe to
# force text entry - exactly like a user typing '584.66 into the cell.
# That quote-prefix also nudges Excel to allocate a new cell style, so
# we explicitly reset each touched cell back to the "Normal" style,
# matching the unstyled cells already on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> @{ column = newValue }
$updates = [ordered]@{
    2  = @{ D = '63.292.13';  E = '  +2.70%  ' }
    3  = @{ D = '3.485.98';   E = '  +2.72%  ' }
    4  = @{ E = '  +0.06%  ' }
    5  = @{ D = '584.66';     E = '  +1.20%  ' }
    6  = @{ D = '147.98';     E = '  +5.20%  ' }
    7  = @{ E = '  +0.00%  ' }
    8  = @{ D = '0.479';      E = '  +0.77%  ' }
    9  = @{ D = '7.70';       E = '  +0.55%  ' }
    10 = @{ E = '  +2.92%  ' }
    11 = @{ E = '  +2.64%  ' }
    12 = @{ D = '4.084.29';   E = '  +2.86%  ' }
    13 = @{ E = '  +4.34%  ' }
    14 = @{ E = '  -0.26%  ' }
    15 = @{ D = '3.486.43';   E = '  +2.85%  ' }
    16 = @{ E = '  +2.03%  ' }
    17 = @{ D = '63.319.97';  E = '  +2.88%  ' }
    18 = @{ D = '6.33';       E = '  +2.40%  ' }
    19 = @{ E = '  +5.39%  ' }
    20 = @{ D = '9.36';       E = '  +4.19%  ' }
    21 = @{ D = '390.88';     E = '  +0.32%  ' }
    22 = @{ E = '  +1.54%  ' }
    23 = @{ D = '75.20';      E = '  +0.10%  ' }
    24 = @{ E = '  -0.06%  ' }
    25 = @{ E = '  +5.21%  ' }
    26 = @{ D = '3.629.45';   E = '  +3.08%  ' }
    27 = @{ E = '  -4.57%  ' }
    28 = @{ D = '7.85';       E = '  +8.02%  ' }
    29 = @{ E = '  -0.04%  ' }
    30 = @{ E = '  +3.11%  ' }
    31 = @{ D = '1.47';       E = '  +6.98%  ' }
    32 = @{ E = '  +0.11%  ' }
    33 = @{ E = '  -0.03%  ' }
    34 = @{ D = '23.84';      E = '  +1.61%  ' }
    35 = @{ D = '5.36';       E = '  +6.08%  ' }
    36 = @{ B = 'Aptos';       C = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt';       D = '7.15';  E = '  +2.77%  ' }
    37 = @{ B = 'EnergySwap';  C = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens';   D = '31.93'; E = '  +23.35%  ' }
    38 = @{ D = '171.36';     E = '  +2.15%  ' }
    39 = @{ D = '1.57';       E = '  +6.78%  ' }
    40 = @{ D = '3.523.41';   E = '  +2.82%  ' }
    41 = @{ D = '0.0770';     E = '  +0.18%  ' }
    42 = @{ E = '  +3.72%  ' }
    43 = @{ E = '  +1.70%  ' }
    44 = @{ E = '  +4.35%  ' }
    45 = @{ D = '42.47';      E = '  +0.02%  ' }
    46 = @{ D = '1.21';       E = '  +6.92%  ' }
    47 = @{ D = '2.621.81';   E = '  +6.63%  ' }
    48 = @{ D = '23.71';      E = '  +4.74%  ' }
    49 = @{ E = '  +13.23%  ' }
    50 = @{ E = '  +1.20%  ' }
    51 = @{ E = '  +3.50%  ' }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $value = $cols[$col]

        $looksNumeric = $value -match '^[+-]?\d+(\.\d+)?$'

        $cell = $ws.Range("$col$row")
        if ($looksNumeric) {
            # Force text entry so "7.70" / "0.0770" keep their exact
            # printed form instead of becoming the numbers 7.7 / 0.077.
            $cell.Value = "'" + $value
        } else {
            $cell.Value = $value
        }
        # Undo the quote-prefix/number-format style Excel assigns when a
        # text override is applied, so the cell stays on the sheet's
        # default (unstyled) look, same as before the edit.
        $cell.Style = 'Normal'
    }
}
